# Insert a new weekly price record at row 203 (Hortaliza / Perejil,
# Terminal La Palmera de La Serena), pushing the existing rows 203:212
# down to 204:213.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 203:212 down one row, carrying formatting (e.g. the date
# style on column D) along with them, just like Excel's own
# Rows.Insert does.
$ws.Rows.Item(203).EntireRow.Insert()

# Populate the newly opened row with this week's data.
$ws.Cells.Item(203, 1).Value = 8
$ws.Cells.Item(203, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(203, 3).Value = "Coquimbo"
$ws.Cells.Item(203, 4).Value = 45008
$ws.Cells.Item(203, 5).Value = 4
$ws.Cells.Item(203, 6).Value = 100112044
$ws.Cells.Item(203, 7).Value = "Perejil"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 2300
$ws.Cells.Item(203, 11).Value = 1800
$ws.Cells.Item(203, 12).Value = 2000
$ws.Cells.Item(203, 13).Value = 1900
$ws.Cells.Item(203, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(203, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(203, 16).Value = 1267
$ws.Cells.Item(203, 17).Value = 1.5
$ws.Cells.Item(203, 18).Value = "Hortaliza"
